# Generate Report for Handoff
#
# The localization pipeline re-ran and produced a fresh handoff: the status
# that used to read "In Translation" is now "Ready for handoff", and the
# handoff/generate timestamps were bumped to the moment the new report was
# produced. Reflect that on all three sheets (Overview + the two per-locale
# detail sheets) and widen the now-longer "Status" columns to fit the text,
# same as the report generator does.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: zh-cn / de-de status cells + HO Xliff generate date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-13 12:44:44"

# --- zh-cn detail sheet: Status + Latest Handoff Datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-13 12:44:35"

# --- de-de detail sheet: Status + Latest Handoff Datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-13 12:44:44"

# --- Widen the Status columns so the longer "Ready for handoff" text fits ---
# (Excel snaps ColumnWidth to whole-pixel increments, so feed it the
# character-width bucket that rounds to the generator's target width.)
$newStatusWidth = 16.3333333333333
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth
